$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for week 41 and 42 (semana 42 de 2025 update)
$ws.Range("B42").Value = 348
$ws.Range("B43").Value = 317

# Add new row for week 43
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 56
